# Offline attendance system update: add new check-in/out records (Linh, bb)
# to the Attendance log, and refresh the Monthly Report summary (sorted by
# name) with the new entries inserted in their alphabetical position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Attendance  -- append the two new raw attendance rows
# ---------------------------------------------------------------------
$att = $wb.Worksheets.Item("Attendance")

# Row 10: Linh, clocked in late, no clock-out recorded yet
$att.Cells.Item(10, 1).Value = "Linh"
# Force the "Day" column to stay plain text -- "2026-01-10" looks like an
# ISO date and Excel would otherwise silently convert it to a date serial
# number. Mark the cell as Text, type the value, then drop back to the
# sheet's default (unstyled) look so formatting matches the other rows.
$att.Cells.Item(10, 2).NumberFormat = "@"
$att.Cells.Item(10, 2).Value = "2026-01-10"
$att.Cells.Item(10, 2).Style = "Normal"
$att.Cells.Item(10, 3).Value = "09:18:51"
$att.Cells.Item(10, 6).Value = "LATE"
$att.Cells.Item(10, 7).Value = "1h 18m"
$att.Cells.Item(10, 8).Value = "0m"

# Carry over the "LATE" status formatting (bold red, style index 2) from
# an existing LATE row so the new cells match exactly.
$att.Range("F2:G2").Copy()
$att.Range("F10:G10").PasteSpecial(-4122)

# Row 11: bb, clocked in late and already clocked out
$att.Cells.Item(11, 1).Value = "bb"
$att.Cells.Item(11, 2).NumberFormat = "@"
$att.Cells.Item(11, 2).Value = "2026-01-10"
$att.Cells.Item(11, 2).Style = "Normal"
$att.Cells.Item(11, 3).Value = "11:24:40"
$att.Cells.Item(11, 4).Value = "11:34:34"
$att.Cells.Item(11, 5).Value = "0h 9m"
$att.Cells.Item(11, 6).Value = "LATE"
$att.Cells.Item(11, 7).Value = "3h 24m"
$att.Cells.Item(11, 8).Value = "0m"

$att.Range("F2:G2").Copy()
$att.Range("F11:G11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet 2: Monthly Report -- insert the two new summary rows in their
# alphabetically-sorted position (report is kept sorted by Name).
# ---------------------------------------------------------------------
$rep = $wb.Worksheets.Item("Monthly Report")

# "Linh" sorts between "FinalTest" (row 3) and "TestUser2" (row 4) --
# insert a new row 4 and push the rest down.
$rep.Rows.Item(4).Insert()
$rep.Cells.Item(4, 1).Value = "Linh"
$rep.Cells.Item(4, 2).Value = 1
$rep.Cells.Item(4, 3).Value = 0
$rep.Cells.Item(4, 4).Value = "1h 18m"
$rep.Cells.Item(4, 5).Value = "0m"

# "bb" sorts between "User1" (now row 7) and "gg" (now row 8) -- insert
# a new row 8 and push "gg" down to row 9.
$rep.Rows.Item(8).Insert()
$rep.Cells.Item(8, 1).Value = "bb"
$rep.Cells.Item(8, 2).Value = 1
$rep.Cells.Item(8, 3).Value = 0
$rep.Cells.Item(8, 4).Value = "3h 24m"
$rep.Cells.Item(8, 5).Value = "0m"
